# Rename wc_lang DfbaNetReaction -> DfbaObjReaction, DfbaNetSpecies -> DfbaObjSpecies
# Reflects the sheet/label renames in the workbook:
#   "dFBA net reactions" -> "dFBA objective reactions"
#   "dFBA net species"   -> "dFBA objective species"
# and the matching header text "dFBA net reaction" -> "dFBA objective reaction".

$wb = $excel.ActiveWorkbook

# 1) Rename the two worksheets. Excel automatically keeps every defined name
#    (_xlnm._FilterDatabase, _FilterDatabase_0, etc.) that references these
#    sheets in sync with the new names.
$wsReactions = $wb.Worksheets.Item("dFBA net reactions")
$wsReactions.Name = "dFBA objective reactions"

$wsSpecies = $wb.Worksheets.Item("dFBA net species")
$wsSpecies.Name = "dFBA objective species"

# 2) Update the header cell text that mirrored the old sheet label.
$wsSpecies.Range("C1").Value = "dFBA objective reaction"

# 3) Make "dFBA objective species" the active sheet/tab, with F7 selected,
#    matching the saved view state in the workbook.
$wsSpecies.Activate() | Out-Null
$wsSpecies.Range("F7").Select() | Out-Null
